$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-10: qL (S), qR (T), dload (U) become a successful line of thrust load
$ws.Range("S8:T10").Value = 3000
$ws.Range("U8:U10").Value = 38107.2171642066

# Rows 11-12: previously loaded, now cleared to zero
$ws.Range("S11:U12").Value = 0
